$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 26, shifting existing rows (26-122) down to (27-123)
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with data
$ws.Cells.Item(26, 1).Value = 9
$ws.Cells.Item(26, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(26, 3).Value = "Metropolitana"
$ws.Cells.Item(26, 4).Value = 45063
$ws.Cells.Item(26, 5).Value = 13
$ws.Cells.Item(26, 6).Value = 100112005
$ws.Cells.Item(26, 7).Value = "Puerro"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 70
$ws.Cells.Item(26, 11).Value = 7000
$ws.Cells.Item(26, 12).Value = 7000
$ws.Cells.Item(26, 13).Value = 7000
$ws.Cells.Item(26, 14).Value = "`$/paquete 20 unidades"
$ws.Cells.Item(26, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(26, 16).Value = 350
$ws.Cells.Item(26, 17).Value = 20
$ws.Cells.Item(26, 18).Value = "Hortaliza"
